# [Kadastro App] Yeni kayit eklendi: 2936
#
# Appends a new record row to both the "Kayitlar" master sheet and the
# per-district "Erdemli" sheet. The source workbook stores every cell
# (including numeric-looking values like "2936" or "1", and the date
# "2025-09-09") as literal TEXT rather than numbers/dates, so a plain
# Range.Value assignment (which lets Excel auto-detect numbers/dates)
# would change the stored type. To keep the new cells as genuine text
# -- the same way the rest of the sheet is encoded -- each value is
# written as a `="..."` text-formula and then converted to a plain
# value in place via Copy / PasteSpecial (values only), which drops the
# formula but keeps the string result and the default cell style.

$wb = $excel.ActiveWorkbook

$record = @("2936", "2025-09-09", "Erdemli", "1", "ÇAP", "SEVİL SARAÇER (Tekniker)")
$cols = @("A", "B", "C", "D", "E", "F")

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    for ($i = 0; $i -lt $record.Length; $i++) {
        $cellRef = "$($cols[$i])$newRow"
        $escaped = $record[$i].Replace('"', '""')
        $ws.Range($cellRef).Formula = "=""$escaped"""
    }

    $rowRange = $ws.Range("A" + $newRow + ":F" + $newRow)
    $rowRange.Copy()
    $rowRange.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
